$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows per the diff
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -2
$ws.Range("F7").Value = -13
$ws.Range("F9").Value = 1
$ws.Range("F11").Value = -3
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("F18").Value = -2
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = -1
